$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3816.8
$ws.Range("J17").Value = 3086.6155
$ws.Range("L17").Value = 9259.8465
$ws.Range("N17").Value = -9595.8465

$ws.Range("H28").Value = 195.13333
$ws.Range("I28").Value = 140.41667
$ws.Range("J28").Value = 414
$ws.Range("K28").Value = 140.41667
$ws.Range("L28").Value = 414
$ws.Range("M28").Value = 344.58333
$ws.Range("N28").Value = -1384

$ws.Range("H137").Value = 2911.3076
$ws.Range("I137").Value = 2760.7778
$ws.Range("K137").Value = 8282.3334
$ws.Range("M137").Value = -5732.3334

$ws.Range("H138").Value = 2336.554
$ws.Range("I138").Value = 2440.7188
$ws.Range("J138").Value = 2235.5454
$ws.Range("K138").Value = 7322.1564
$ws.Range("L138").Value = 6706.6362
$ws.Range("M138").Value = -2182.1564
$ws.Range("N138").Value = -16986.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4119.7837
$ws.Range("I32").Value = 2576.0598
$ws.Range("J32").Value = 18895.428
$ws.Range("K32").Value = 2576.0598
$ws.Range("L32").Value = 18895.428
$ws.Range("M32").Value = -2289.0598
$ws.Range("N32").Value = -19469.428

$ws.Range("H45").Value = 6001440
$ws.Range("I45").Value = 15000818
$ws.Range("K45").Value = 15000818
$ws.Range("M45").Value = -15000441

$ws.Range("H63").Value = 1501.6666
$ws.Range("I63").Value = 1422
$ws.Range("K63").Value = 1422
$ws.Range("M63").Value = -736

$ws.Range("H66").Value = 1501.6666
$ws.Range("I66").Value = 1422
$ws.Range("K66").Value = 7110
$ws.Range("M66").Value = -3678

$ws.Range("H74").Value = 1301.6154
$ws.Range("I74").Value = 491.17648
$ws.Range("J74").Value = 2832.4443
$ws.Range("K74").Value = 491.17648
$ws.Range("L74").Value = 2832.4443
$ws.Range("M74").Value = 382.82352
$ws.Range("N74").Value = -4580.4443

$ws.Range("H77").Value = 1301.6154
$ws.Range("I77").Value = 491.17648
$ws.Range("J77").Value = 2832.4443
$ws.Range("K77").Value = 2455.8824
$ws.Range("L77").Value = 14162.2215
$ws.Range("M77").Value = 1912.1176
$ws.Range("N77").Value = -22898.2215

$ws.Range("H135").Value = 31579.6
$ws.Range("J135").Value = 31579.6
$ws.Range("L135").Value = 31579.6
$ws.Range("N135").Value = -41719.6

$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 19300

$ws.Range("H85").Value = 19300

$ws.Range("H107").Value = 912.4666999999999
$ws.Range("I107").Value = 656.3333
$ws.Range("J107").Value = 1083.2222
$ws.Range("K107").Value = 656.3333
$ws.Range("L107").Value = 1083.2222
$ws.Range("M107").Value = 1263.6667
$ws.Range("N107").Value = -4923.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 700.6
$ws.Range("I16").Value = 743.5714
$ws.Range("J16").Value = 600.3333
$ws.Range("K16").Value = 743.5714
$ws.Range("L16").Value = 600.3333
$ws.Range("M16").Value = -456.5714
$ws.Range("N16").Value = -1174.3333

$ws.Range("H31").Value = 2161.8572
$ws.Range("I31").Value = 1857
$ws.Range("J31").Value = 2314.2856
$ws.Range("K31").Value = 1857
$ws.Range("L31").Value = 2314.2856
$ws.Range("M31").Value = -1562
$ws.Range("N31").Value = -2904.2856

$ws.Range("H34").Value = 2161.8572
$ws.Range("I34").Value = 1857
$ws.Range("J34").Value = 2314.2856
$ws.Range("K34").Value = 1857
$ws.Range("L34").Value = 2314.2856
$ws.Range("M34").Value = -1655
$ws.Range("N34").Value = -2718.2856

$ws.Range("H113").Value = 700.6
$ws.Range("I113").Value = 743.5714
$ws.Range("J113").Value = 600.3333
$ws.Range("K113").Value = 743.5714
$ws.Range("L113").Value = 600.3333
$ws.Range("M113").Value = 1426.4286
$ws.Range("N113").Value = -4940.3333

$ws.Range("H122").Value = 3116.7144
$ws.Range("I122").Value = 2051.75
$ws.Range("K122").Value = 6155.25
$ws.Range("M122").Value = -3705.25

$ws.Range("H132").Value = 2365.5217
$ws.Range("I132").Value = 1453.2667
$ws.Range("K132").Value = 4359.800099999999
$ws.Range("M132").Value = -1829.800099999999

$ws.Range("H134").Value = 1679.5518
$ws.Range("I134").Value = 1323.24
$ws.Range("K134").Value = 3969.72
$ws.Range("M134").Value = -1434.72

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 16220.8
$ws.Range("I17").Value = 98
$ws.Range("J17").Value = 20251.5
$ws.Range("K17").Value = 294
$ws.Range("L17").Value = 60754.5
$ws.Range("M17").Value = -125
$ws.Range("N17").Value = -61092.5

$ws.Range("H131").Value = 11417.8
$ws.Range("J131").Value = 12027.817
$ws.Range("L131").Value = 36083.451
$ws.Range("N131").Value = -46163.451

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 21563.143
$ws.Range("J26").Value = 21563.143
$ws.Range("L26").Value = 21563.143
$ws.Range("N26").Value = -22123.143

$ws.Range("H50").Value = 21563.143
$ws.Range("J50").Value = 21563.143
$ws.Range("L50").Value = 21563.143
$ws.Range("N50").Value = -22559.143

$ws.Range("H132").Value = 2567068.8
$ws.Range("I132").Value = 3848013.5
$ws.Range("K132").Value = 11544040.5
$ws.Range("M132").Value = -11541510.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8318
$ws.Range("I16").Value = 10840.667
$ws.Range("J16").Value = 750
$ws.Range("K16").Value = 10840.667
$ws.Range("L16").Value = 750
$ws.Range("M16").Value = -10670.667
$ws.Range("N16").Value = -1090

$ws.Range("H25").Value = 16000
$ws.Range("I25").Value = 15000
$ws.Range("J25").Value = 16500
$ws.Range("K25").Value = 15000
$ws.Range("L25").Value = 16500
$ws.Range("M25").Value = -14770
$ws.Range("N25").Value = -16960

$ws.Range("H46").Value = 1887.3846
$ws.Range("J46").Value = 2794
$ws.Range("L46").Value = 2794
$ws.Range("N46").Value = -3170

$ws.Range("H132").Value = 1858.091
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1858.091
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 5574.272999999999
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -10634.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H46").Value = 29998
$ws.Range("J46").Value = 29998
$ws.Range("L46").Value = 29998
$ws.Range("N46").Value = -30460

$ws.Range("H134").Value = 29998
$ws.Range("J134").Value = 29998
$ws.Range("L134").Value = 89994
$ws.Range("N134").Value = -95064
